$wb = $excel.ActiveWorkbook

# --- Emissions sheet: insert "emission_name" column and update values ---
$emissions = $wb.Worksheets.Item("Emissions")

# Move former column B ("emission_unit" header, "ton" value) to column C,
# copying its header formatting (bold, bordered, centered/top) as well.
$emissions.Range("C1").Value = "emission_unit"
$emissions.Range("C1").Font.Bold = $true
$emissions.Range("C1").HorizontalAlignment = -4108
$emissions.Range("C1").VerticalAlignment = -4160
$emissions.Range("C1").Borders.LineStyle = 1

$emissions.Range("C2").Value = "ton"

# New column B: emission_name
$emissions.Range("B1").Value = "emission_name"
$emissions.Range("B1").Font.Bold = $true
$emissions.Range("B1").HorizontalAlignment = -4108
$emissions.Range("B1").VerticalAlignment = -4160
$emissions.Range("B1").Borders.LineStyle = 1

$emissions.Range("B2").Value = "CO2 emissions"

# Column A: emission_name default changes from "CO2-equivalent" to "CO2"
$emissions.Range("A2").Value = "CO2"

$emissions.Range("C10").Select()

# --- Techs sheet: change the selected cell ---
$techs = $wb.Worksheets.Item("Techs")
$techs.Range("C2").Select()

# --- Workbook window position ---
$excel.Windows.Item(1).Left = 7300
$excel.Windows.Item(1).Top = 700

$wb.Save()
